$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1/J1 should look like the existing header cells (B1..H1):
# bold font, thin border, centered/top aligned (style index 1 in the sheet).
# Copy the formatting from H1 (an existing header) onto I1:J1, then set the
# text separately so the shared style is reused instead of a new one minted.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I (I0) and J (IF), rows 2-13
$data = @{
    2  = @(8, 8)
    3  = @(12, 13)
    4  = @(8, 8)
    5  = @(6, 8)
    6  = @(7, 7)
    7  = @(7, 7)
    8  = @(7, 7)
    9  = @(6, 8)
    10 = @(1, 6)
    11 = @(1, 6)
    12 = @(1, 3)
    13 = @(1, 3)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
